{"js": "// The document contains a single table of two-digit \u00f7 one-digit division\n// problems. Every populated cell's \"answer\" text is being replaced with a\n// new problem/answer string, in the same left-to-right, top-to-bottom\n// reading order the cells already appear in. Blank filler rows/cells are\n// left untouched.\nconst replacements = [\n  \"19\u00f79=2, 1|51\u00f72=25, 1\",\n  \"87\u00f75=17, 2|17\u00f78=2, 1\",\n  \"31\u00f73=10, 1|90\u00f75=18, 0\",\n  \"59\u00f77=8, 3|41\u00f75=8, 1\",\n  \"67\u00f79=7, 4|74\u00f74=18, 2\",\n  \"95\u00f77=13, 4|36\u00f72=18, 0\",\n  \"93\u00f78=11, 5|60\u00f73=20, 0\",\n  \"39\u00f79=4, 3|15\u00f72=7, 1\",\n  \"47\u00f76=7, 5|68\u00f75=13, 3\",\n  \"94\u00f74=23, 2|76\u00f76=12, 4\",\n  \"35\u00f72=17, 1|71\u00f79=7, 8\",\n  \"77\u00f77=11, 0|62\u00f73=20, 2\",\n  \"16\u00f73=5, 1|66\u00f73=22, 0\",\n  \"17\u00f78=2, 1|84\u00f75=16, 4\",\n  \"23\u00f78=2, 7|40\u00f76=6, 4\",\n  \"79\u00f78=9, 7|48\u00f75=9, 3\",\n  \"65\u00f74=16, 1|25\u00f72=12, 1\",\n  \"20\u00f75=4, 0|92\u00f72=46, 0\",\n  \"39\u00f79=4, 3|97\u00f76=16, 1\",\n  \"39\u00f76=6, 3|53\u00f79=5, 8\",\n  \"84\u00f73=28, 0|61\u00f73=20, 1\",\n  \"91\u00f78=11, 3|87\u00f73=29, 0\",\n  \"64\u00f76=10, 4|12\u00f79=1, 3\",\n  \"56\u00f77=8, 0|92\u00f77=13, 1\",\n  \"49\u00f73=16, 1|95\u00f75=19, 0\",\n].map((s) => {\n  const i = s.indexOf(\"|\");\n  return { oldText: s.slice(0, i), newText: s.slice(i + 1) };\n});\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet next = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const cellText = grid[r][c];\n    if (cellText === \"\" || cellText === undefined || cellText === null) continue;\n    if (next >= replacements.length) continue;\n    const { oldText, newText } = replacements[next];\n    next++;\n    if (cellText !== oldText) {\n      throw new Error(\n        `Cell (${r},${c}) text \"${cellText}\" did not match expected \"${oldText}\"`\n      );\n    }\n    table.getCell(r, c).value = newText;\n  }\n}\n\nawait context.sync();\n\nif (next !== replacements.length) {\n  throw new Error(`Only applied ${next} of ${replacements.length} replacements`);\n}\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit division\n# problems. Every populated cell's \"answer\" text is being replaced with a\n# new problem/answer string, in the same left-to-right, top-to-bottom\n# reading order the cells already appear in. Blank filler rows/cells are\n# left untouched.\n$oldTexts = @(\n    \"19\u00f79=2, 1\",\n    \"87\u00f75=17, 2\",\n    \"31\u00f73=10, 1\",\n    \"59\u00f77=8, 3\",\n    \"67\u00f79=7, 4\",\n    \"95\u00f77=13, 4\",\n    \"93\u00f78=11, 5\",\n    \"39\u00f79=4, 3\",\n    \"47\u00f76=7, 5\",\n    \"94\u00f74=23, 2\",\n    \"35\u00f72=17, 1\",\n    \"77\u00f77=11, 0\",\n    \"16\u00f73=5, 1\",\n    \"17\u00f78=2, 1\",\n    \"23\u00f78=2, 7\",\n    \"79\u00f78=9, 7\",\n    \"65\u00f74=16, 1\",\n    \"20\u00f75=4, 0\",\n    \"39\u00f79=4, 3\",\n    \"39\u00f76=6, 3\",\n    \"84\u00f73=28, 0\",\n    \"91\u00f78=11, 3\",\n    \"64\u00f76=10, 4\",\n    \"56\u00f77=8, 0\",\n    \"49\u00f73=16, 1\"\n)\n$newTexts = @(\n    \"51\u00f72=25, 1\",\n    \"17\u00f78=2, 1\",\n    \"90\u00f75=18, 0\",\n    \"41\u00f75=8, 1\",\n    \"74\u00f74=18, 2\",\n    \"36\u00f72=18, 0\",\n    \"60\u00f73=20, 0\",\n    \"15\u00f72=7, 1\",\n    \"68\u00f75=13, 3\",\n    \"76\u00f76=12, 4\",\n    \"71\u00f79=7, 8\",\n    \"62\u00f73=20, 2\",\n    \"66\u00f73=22, 0\",\n    \"84\u00f75=16, 4\",\n    \"40\u00f76=6, 4\",\n    \"48\u00f75=9, 3\",\n    \"25\u00f72=12, 1\",\n    \"92\u00f72=46, 0\",\n    \"97\u00f76=16, 1\",\n    \"53\u00f79=5, 8\",\n    \"61\u00f73=20, 1\",\n    \"87\u00f73=29, 0\",\n    \"12\u00f79=1, 3\",\n    \"92\u00f77=13, 1\",\n    \"95\u00f75=19, 0\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$next = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        # Cell text includes the trailing end-of-cell marker(s); strip any\n        # trailing control characters to get the visible content.\n        $raw = $cellRange.Text\n        $visible = $raw.TrimEnd([char]13, [char]7)\n        if ($visible.Length -eq 0) {\n            continue\n        }\n        if ($next -ge $oldTexts.Length) {\n            continue\n        }\n        $expected = $oldTexts[$next]\n        if ($visible -ne $expected) {\n            throw \"Cell ($r,$c) text '$visible' did not match expected '$expected'\"\n        }\n        $cellRange.Text = $newTexts[$next]\n        $next = $next + 1\n    }\n}\n\nif ($next -ne $oldTexts.Length) {\n    throw \"Only applied $next of $($oldTexts.Length) replacements\"\n}\n"}
